# Fix the nutrient-quantity table: a prior step inserted the 3 new nutrient rows
# (AGPOLI, AGTRANS, AGSAT) without re-aligning the food-quantity columns, so every
# row below FIBRA was reading the PREVIOUS (now-stale/empty) food quantity instead
# of being recomputed for the new, larger nutrient list. This reassigns the correct
# nutrient label + quantities to every shifted row and appends the 3 rows (FOLATO,
# FOSFORO, ZINCO) that the old, too-short table was missing at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 502.35

$ws.Range("D3").Value = 96.88

$ws.Range("D4").Value = 24.12
$ws.Range("E4").Value = 238.884

$ws.Range("D5").Value = 7.64
$ws.Range("E5").Value = 358.326

$ws.Range("D6").Value = 25.26

$ws.Range("B7").Value = "AGPOLI"
$ws.Range("C7").Value = 11.47
$ws.Range("D7").Value = 3.87

$ws.Range("B8").Value = "AGTRANS"
$ws.Range("C8").Value = 0.9
$ws.Range("D8").Value = 0.17

$ws.Range("B9").Value = "AGSAT"
$ws.Range("C9").Value = 11.01
$ws.Range("D9").Value = 1.53
$ws.Range("E9").Value = 23.8884

$ws.Range("B10").Value = "COLEST"
$ws.Range("C10").Value = 122.22
$ws.Range("D10").Value = 5.38
$ws.Range("E10").Value = 238.884

$ws.Range("B11").Value = "CALCIO"
$ws.Range("C11").Value = 219.6
$ws.Range("D11").Value = 107.01
$ws.Range("E11").Value = 143.3304

$ws.Range("B12").Value = "SODIO"
$ws.Range("C12").Value = 2623.22
$ws.Range("D12").Value = 1047.94
$ws.Range("E12").Value = 1

$ws.Range("B13").Value = "POTASSIO"
$ws.Range("C13").Value = 2241.43
$ws.Range("D13").Value = 1085.02
$ws.Range("E13").Value = 3510

$ws.Range("B14").Value = "FERRO"
$ws.Range("C14").Value = 13.55
$ws.Range("D14").Value = 7.1
$ws.Range("E14").Value = 6.8

$ws.Range("B15").Value = "MAGNESIO"
$ws.Range("C15").Value = 309.36
$ws.Range("D15").Value = 170.19
$ws.Range("E15").Value = 303

$ws.Range("B16").Value = "TIAMINA"
$ws.Range("C16").Value = 0.7
$ws.Range("D16").Value = 0.36
$ws.Range("E16").Value = 0.9

$ws.Range("B17").Value = "RIBOFLAVINA"
$ws.Range("C17").Value = 0.95
$ws.Range("D17").Value = 0.36
$ws.Range("E17").Value = 1

$ws.Range("B18").Value = "NIACINA"
$ws.Range("C18").Value = 11.74
$ws.Range("D18").Value = 1.65
$ws.Range("E18").Value = 11.5

$ws.Range("B19").Value = "PIRIDOXAMINA"
$ws.Range("C19").Value = 0.52
$ws.Range("D19").Value = 0.35
$ws.Range("E19").Value = 1.1

$ws.Range("B20").Value = "COBALAMINA"
$ws.Range("C20").Value = 4.06
$ws.Range("D20").Value = 0.18
$ws.Range("E20").Value = 2

$ws.Range("B21").Value = "VITC"
$ws.Range("C21").Value = 3.74
$ws.Range("D21").Value = 1.62
$ws.Range("E21").Value = 66.09999999999999

$ws.Range("B22").Value = "VITA_RAE"
$ws.Range("C22").Value = 41.23
$ws.Range("D22").Value = 1.61
$ws.Range("E22").Value = 560

$ws.Range("B23").Value = "COBRE"
$ws.Range("C23").Value = 1.39
$ws.Range("D23").Value = 0.78
$ws.Range("E23").Value = 0.7

# Append the 3 new rows (24-26); copy row 23 formatting/style first, then set values
$ws.Range("A23:E23").Copy() | Out-Null
$ws.Range("A24:E26").PasteSpecial(-4122) | Out-Null

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "FOLATO"
$ws.Range("C24").Value = 481.27
$ws.Range("D24").Value = 406.45
$ws.Range("E24").Value = 322

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "FOSFORO"
$ws.Range("C25").Value = 989.78
$ws.Range("D25").Value = 396.61
$ws.Range("E25").Value = 8

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "ZINCO"
$ws.Range("C26").Value = 16.71
$ws.Range("D26").Value = 3.57
$ws.Range("E26").Value = 649